$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: semantic type annotations for Mes nombre (D), Grupo de tipo de jornada (E),
# Sexo (K) and Mes y año (L) switch from "measure" to "dimension".
$ws.Range("D3").Value = "iaest-dimension:mes-nombre"
$ws.Range("E3").Value = "iaest-dimension:grupo-de-tipo-de-jornada"
$ws.Range("K3").Value = "iaest-dimension:sexo"
$ws.Range("L3").Value = "iaest-dimension:mes-y-ano"

# Row 4: dim/medida marker for those same columns switches from "medida" to "dim".
$ws.Range("D4").Value = "dim"
$ws.Range("E4").Value = "dim"
$ws.Range("K4").Value = "dim"
$ws.Range("L4").Value = "dim"

# Row 5: data type for the now-dimension columns becomes "skos:Concept"
# (Mes y año keeps its scalar "xsd:string" type).
$ws.Range("D5").Value = "skos:Concept"
$ws.Range("E5").Value = "skos:Concept"
$ws.Range("K5").Value = "skos:Concept"
$ws.Range("L5").Value = "xsd:string"

# New row 6: mapping files for the new dimension columns, formatted like the
# rest of the table (copy the style down from row 5).
$ws.Range("D6").Value = "mapping-mes-nombre.xlsx"
$ws.Range("E6").Value = "mapping-grupo-de-tipo-de-jornada.xlsx"
$ws.Range("K6").Value = "mapping-sexo.xlsx"

$ws.Range("D5:E5").Copy()
$ws.Range("D6:E6").PasteSpecial(-4122)
$ws.Range("K5").Copy()
$ws.Range("K6").PasteSpecial(-4122)
